# Trade #68 closed at 2026-02-17 12:54:00 - unknown UNKNOWN +0.000%
#
# Updates:
#  - Summary sheet: Total Trades 67 -> 68, Win Rate % 44.78 -> 44.12
#  - Strategy Status sheet: MarketMaking row Trades 67 -> 68, Win Rate % 44.78 -> 44.12
#  - All Trades sheet: append new trade row (row 69)
#  - MarketMaking sheet: append the same new trade row (row 69)

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B6").Value = 68
$summary.Range("B9").Value = 44.12

# ---------------------------------------------------------------------------
# Strategy Status sheet (MarketMaking is row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("D4").Value = 68
$status.Range("G4").Value = 44.12

# ---------------------------------------------------------------------------
# Helper: write the new closed-trade row (row 69) onto a trades-log sheet
# ---------------------------------------------------------------------------
function Add-Trade68Row($ws) {
    $ws.Range("A69").Value = 68

    # Date/Time columns look like dates/times so they must be forced to
    # plain text, otherwise Excel auto-converts them into date/time serials.
    $ws.Range("B69").NumberFormat = "@"
    $ws.Range("B69").Value = "2026-02-17"
    $ws.Range("C69").NumberFormat = "@"
    $ws.Range("C69").Value = "12:53:53"
    $ws.Range("B69:C69").ClearFormats()

    $ws.Range("D69").Value = "MarketMaking"
    $ws.Range("E69").Value = "DOWN"
    $ws.Range("F69").Value = 0.93
    $ws.Range("G69").Value = 0.93
    $ws.Range("H69").Value = "CLOSED"
    $ws.Range("I69").Value = 0
    $ws.Range("J69").Value = 0
    $ws.Range("K69").Value = 100.16
    $ws.Range("L69").Value = 0
    $ws.Range("M69").Value = 0
    $ws.Range("N69").Value = 0.6
    $ws.Range("O69").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P69").Value = "early_exit"
    $ws.Range("Q69").Value = 0.13
}

# ---------------------------------------------------------------------------
# All Trades sheet
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")
Add-Trade68Row $allTrades

# ---------------------------------------------------------------------------
# MarketMaking sheet (same trade log, filtered to this strategy)
# ---------------------------------------------------------------------------
$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-Trade68Row $marketMaking
